$wb = $excel.ActiveWorkbook

# The status text "Ready for handoff" is shared (as one shared string) across
# the Overview sheet (columns for zh-cn/de-de on the 502d1cb6 file row) and the
# per-locale "Status" column on the zh-cn/de-de sheets. Updating every cell
# that currently holds that text keeps them all pointing at the same new
# string, matching how the shared string table collapses identical text.
$newStatus = "Handback transform failed"
$oldStatus = "Ready for handoff"

foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    foreach ($cell in $used.Cells) {
        # Keep the string literal on the left of -eq: some cells hold a
        # Boolean (True/False), and Boolean -eq String coerces the string to
        # Boolean (any non-empty text -> $true), causing false matches.
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value = $newStatus
        }
    }
}

# The worksheet's stored <col width> is ColumnWidth plus a fixed ~5/6
# character padding (standard Excel "characters -> XML width" conversion).
# Using 40 directly would persist as ~40.83; back the padding out so the
# saved width lands on an even 40 (matches the wider "Error Detail" column
# now that it holds long diagnostic text).
$errorDetailColumnWidth = 40 - (5 / 6)

# zh-cn sheet: record the handback-transform error detail for the
# 502d1cb6 file row (row 3) in column P ("Error Detail").
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("P3").Value = "Handback file name: 2v442sxe.cms is different with handoff file name: 502d1cb6-495d-4a80-8648-e5970129dfa0.01356ca608a0cfed77ead40f326e06d30264563c.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = $errorDetailColumnWidth

# de-de sheet: same, for the de-de handoff artifact.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("P3").Value = "Handback file name: 2v442sxe.cms is different with handoff file name: 502d1cb6-495d-4a80-8648-e5970129dfa0.01356ca608a0cfed77ead40f326e06d30264563c.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = $errorDetailColumnWidth
